$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert a new row at 251 (pushes old rows 251-347 down to 252-348)
$ws.Rows.Item(251).Insert()

# Step 2: insert 3 new rows at 327 (after step 1, old row 326 now sits at 327;
# this pushes it - and everything below - down to 330-351)
$ws.Rows.Item(327).Insert()
$ws.Rows.Item(327).Insert()
$ws.Rows.Item(327).Insert()

# Fill in the brand-new row 251
$ws.Range("A251").Value = 5
$ws.Range("B251").Value = "Macroferia Regional de Talca"
$ws.Range("C251").Value = "Maule"
$ws.Range("D251").Value = 45007
$ws.Range("E251").Value = 7
$ws.Range("F251").Value = 100112021
$ws.Range("G251").Value = "Ají"
$ws.Range("H251").Value = "Cristal"
$ws.Range("I251").Value = "Primera"
$ws.Range("J251").Value = 200
$ws.Range("K251").Value = 12000
$ws.Range("L251").Value = 12000
$ws.Range("M251").Value = 12000
$ws.Range("N251").Value = "$/saco 25 kilos"
$ws.Range("O251").Value = "Región del Maule"
$ws.Range("P251").Value = 480
$ws.Range("Q251").Value = 25
$ws.Range("R251").Value = "Hortaliza"

# Fill in the brand-new row 327
$ws.Range("A327").Value = 5
$ws.Range("B327").Value = "Macroferia Regional de Talca"
$ws.Range("C327").Value = "Maule"
$ws.Range("D327").Value = 45008
$ws.Range("E327").Value = 7
$ws.Range("F327").Value = 100112021
$ws.Range("G327").Value = "Ají"
$ws.Range("H327").Value = "Cacho cabra rojo"
$ws.Range("I327").Value = "Primera"
$ws.Range("J327").Value = 20
$ws.Range("K327").Value = 18000
$ws.Range("L327").Value = 18000
$ws.Range("M327").Value = 18000
$ws.Range("N327").Value = "$/saco 25 kilos"
$ws.Range("O327").Value = "Región del Maule"
$ws.Range("P327").Value = 720
$ws.Range("Q327").Value = 25
$ws.Range("R327").Value = "Hortaliza"

# Fill in the brand-new row 328
$ws.Range("A328").Value = 5
$ws.Range("B328").Value = "Macroferia Regional de Talca"
$ws.Range("C328").Value = "Maule"
$ws.Range("D328").Value = 45008
$ws.Range("E328").Value = 7
$ws.Range("F328").Value = 100112021
$ws.Range("G328").Value = "Ají"
$ws.Range("H328").Value = "Cacho cabra verde"
$ws.Range("I328").Value = "Primera"
$ws.Range("J328").Value = 150
$ws.Range("K328").Value = 10000
$ws.Range("L328").Value = 10000
$ws.Range("M328").Value = 10000
$ws.Range("N328").Value = "$/saco 25 kilos"
$ws.Range("O328").Value = "Región del Maule"
$ws.Range("P328").Value = 400
$ws.Range("Q328").Value = 25
$ws.Range("R328").Value = "Hortaliza"

# Fill in the brand-new row 329
$ws.Range("A329").Value = 5
$ws.Range("B329").Value = "Macroferia Regional de Talca"
$ws.Range("C329").Value = "Maule"
$ws.Range("D329").Value = 45008
$ws.Range("E329").Value = 7
$ws.Range("F329").Value = 100112021
$ws.Range("G329").Value = "Ají"
$ws.Range("H329").Value = "Cristal"
$ws.Range("I329").Value = "Primera"
$ws.Range("J329").Value = 200
$ws.Range("K329").Value = 12000
$ws.Range("L329").Value = 12000
$ws.Range("M329").Value = 12000
$ws.Range("N329").Value = "$/saco 25 kilos"
$ws.Range("O329").Value = "Región del Maule"
$ws.Range("P329").Value = 480
$ws.Range("Q329").Value = 25
$ws.Range("R329").Value = "Hortaliza"
